$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidateFormulas")
$ws.Activate()

# Insert a new row before row 31 - shifts existing rows 31..41 down to 32..42
$ws.Rows.Item(31).Insert()

# New header cell + three IF() formulas in the freshly inserted row
$ws.Cells.Item(31, 1).Value = "If"
$ws.Cells.Item(31, 1).Font.Bold = $true

$ws.Cells.Item(31, 2).Formula = "=IF(B2>3,B3,B5)"
$ws.Cells.Item(31, 3).Formula = "=IF((B2*B3)*C1<0,(B2*B3)*C1,ABS((B2*B3)*C1))"
$ws.Cells.Item(31, 4).Formula = "=IF((B2*B3)*C1<0,ABS((B2*B3)*C1),(B2*B3)*C1)"

$ws.Range("D31").Select()
